$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.110.95'
$ws.Range("E2").Value = '  +2.82%  '
$ws.Range("D3").Value = '2.277.85'
$ws.Range("E3").Value = '  +2.74%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.591'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.90%  '
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("E9").Value = '  +2.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0843'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.07%  '
$ws.Range("E13").Value = '  +1.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.887'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.16%  '
$ws.Range("D15").Value = '2.623.81'
$ws.Range("E15").Value = '  +2.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.64'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.39%  '
$ws.Range("D17").Value = '2.274.02'
$ws.Range("E17").Value = '  +3.07%  '
$ws.Range("D18").Value = '44.016.16'
$ws.Range("E18").Value = '  +2.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.47%  '
$ws.Range("E20").Value = '  +4.01%  '
$ws.Range("E21").Value = '  +2.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.83%  '
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '238.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("E25").Value = '  +4.57%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.25'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +14.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.55'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '163.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0887'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.36%  '
$ws.Range("E34").Value = '  -1.45%  '
$ws.Range("E35").Value = '  +5.62%  '
$ws.Range("E36").Value = '  +2.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.115'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +12.37%  '
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.37%  '
$ws.Range("E40").Value = '  +1.51%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0329'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.24%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.40'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +26.27%  '
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").Value = '1.765.57'
$ws.Range("E44").Value = '  -6.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.209'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '86.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '60.06'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("E50").Value = '  +2.77%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.81%  '
